$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New table header row (row 6)
$ws.Range("A6").Value = "Стартовый вес, кг"
$ws.Range("B6").Value = "Текущий вес, кг"
$ws.Range("C6").Value = "Цель, кг"
$ws.Range("D6").Value = "Прогресс, %"

# New data row (row 7) - values are entered as text (quote-prefixed) so they
# stay as literal "100.0"-style strings instead of being parsed as numbers
$ws.Range("A7").Value = "'100.0"
$ws.Range("B7").Value = "'90.0"
$ws.Range("C7").Value = "'50.0"
$ws.Range("D7").Value = "'20.0"

# Give columns A:D a fixed width of 20 (Excel's ColumnWidth property is
# offset from the stored column width by 5/6 of a character)
$ws.Range("A1:D1").ColumnWidth = 19.166666666666668
